$wb = $excel.ActiveWorkbook

# Sheet "展览" updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1102
$ws1.Range("F12").Value = 558
$ws1.Range("F20").Value = 1178
$ws1.Range("F24").Value = 4647
$ws1.Range("F27").Value = 1612
$ws1.Range("F28").Value = 36
$ws1.Range("F29").Value = 83

# Sheet "演出" updates
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value = 42

# Sheet "全部类型" updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F12").Value = 42
$ws4.Range("F15").Value = 1102
$ws4.Range("F23").Value = 558
$ws4.Range("F33").Value = 1178
$ws4.Range("F37").Value = 4647
$ws4.Range("F40").Value = 1612
$ws4.Range("F43").Value = 36
$ws4.Range("F44").Value = 83
